# Work on the "Samples" sheet (the one selected/active in the source file)
# and make it active again so the saved file keeps it as the selected tab.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")
$ws.Activate()

# New column C: header "optionalstrain" and a single value on row 3,
# mirroring the "seekstrain" column B already present.
$ws.Range("C1").Value = "optionalstrain"
$ws.Range("C3").Value = 29823659

# Match the new column widths from the source edit as closely as this
# engine's pixel-quantized ColumnWidth allows.
$ws.Columns.Item(2).ColumnWidth = 11.3333333333333
$ws.Columns.Item(3).ColumnWidth = 11

# Move the selection, as in the edited workbook.
$ws.Range("D4").Select()
